$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos table
# with the latest scraped values. Price strings that look like plain numbers
# (e.g. "213.79") are written with a leading apostrophe so Excel keeps them
# as text instead of auto-converting them to numeric values, matching how
# they were originally stored in the sheet.
$ws.Range('D2').Value = '26.522.13'
$ws.Range('E2').Value = '  -0.85%  '
$ws.Range('D3').Value = '1.626.44'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''213.79'
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('E9').Value = '  -0.40%  '
$ws.Range('E10').Value = '  -0.88%  '
$ws.Range('D11').Value = '''0.0856'
$ws.Range('E11').Value = '  -0.36%  '
$ws.Range('D12').Value = '1.854.25'
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('D13').Value = '1.637.39'
$ws.Range('E13').Value = '  +0.77%  '
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('E15').Value = '  -0.32%  '
$ws.Range('D16').Value = '''63.97'
$ws.Range('E16').Value = '  -1.72%  '
$ws.Range('D17').Value = '''234.79'
$ws.Range('E17').Value = '  -0.69%  '
$ws.Range('D18').Value = '26.530.56'
$ws.Range('E18').Value = '  -0.83%  '
$ws.Range('E19').Value = '  -0.29%  '
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('E22').Value = '  -2.01%  '
$ws.Range('E23').Value = '  -1.64%  '
$ws.Range('E24').Value = '  +0.26%  '
$ws.Range('D25').Value = '''146.10'
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('D27').Value = '''7.07'
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E28').Value = '  -0.81%  '
$ws.Range('D29').Value = '''15.63'
$ws.Range('E29').Value = '  -0.52%  '
$ws.Range('E30').Value = '  -1.13%  '
$ws.Range('E31').Value = '  -0.40%  '
$ws.Range('D32').Value = '1.524.22'
$ws.Range('E32').Value = '  +3.21%  '
$ws.Range('D33').Value = '''3.26'
$ws.Range('E33').Value = '  +0.14%  '
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('E35').Value = '  +2.64%  '
$ws.Range('E36').Value = '  -0.26%  '
$ws.Range('D37').Value = '''0.570'
$ws.Range('E37').Value = '  -0.57%  '
$ws.Range('E38').Value = '  -0.88%  '
$ws.Range('E39').Value = '  -0.59%  '
$ws.Range('E40').Value = '  -1.70%  '
$ws.Range('D42').Value = '''2.22'
$ws.Range('E42').Value = '  +0.47%  '
$ws.Range('D43').Value = '1.766.22'
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').Value = '''62.66'
$ws.Range('E44').Value = '  +0.76%  '
$ws.Range('E45').Value = '  -0.77%  '
$ws.Range('D46').Value = '''0.910'
$ws.Range('E46').Value = '  -4.78%  '
$ws.Range('D47').Value = '''89.90'
$ws.Range('E47').Value = '  +1.55%  '
$ws.Range('E48').Value = '  +0.41%  '

# Row 49 is now Cronos (previously BabyDogeCoin) and row 50 is now EnergySwap
# (previously Cronos); update coin name, link, price, and volume accordingly.
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '''0.0501'
$ws.Range('E49').Value = '  -0.65%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''7.56'
$ws.Range('E50').Value = '  +0.98%  '
$ws.Range('D51').Value = '''0.0964'
$ws.Range('E51').Value = '  -0.39%  '
